$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B6").Value = "d"
$ws.Range("G6").Value = "camp chess 6"
$ws.Range("P6").Value = "01/01/2023"

$ws.Range("B7").Value = "tutor"
$ws.Range("P7").Value = "10/01/2023"
$ws.Range("Q7").Value = "31/01/2023"
$ws.Range("T7").Value = 5

$ws.Range("Q7").Select()
